# Add a header row (Year / Market cap / Change) above the existing
# Amazon market-cap data table, shifting all existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the very top; everything currently on row 1
# (and below) moves down to make room.
$ws.Rows("1:1").Insert()

# The data row that used to be row 1 is now row 2; its "Change" column
# (C) carries the percentage number format we want the new header cell
# to reuse, so copy that formatting onto C1 before writing the header
# text (keeps the same style index instead of minting a new one).
$ws.Range("C2").Copy($ws.Range("C1"))

# Write the new header labels.
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Market cap"
$ws.Range("C1").Value = "Change"

# Re-select the full (now one-row-taller) table, matching the widened
# used range after the insert.
[void]$ws.Range("A1:C29").Select()
